# Update the "Total Jobs" sensitivity sheet with the new values recalculated
# from the IEV model under exact 1-year shifting with no double-deployment
# capex (values copied from the [2021:2050] run into the [2020:2049] columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Breeze
$ws.Range("B2").Value = 146313.3958052042
$ws.Range("C2").Value = 146827.02811716934
$ws.Range("D2").Value = 148035.20963514346
$ws.Range("E2").Value = 137225.56905248243
$ws.Range("F2").Value = 133646.7174088397
$ws.Range("G2").Value = 123471.41113176254
$ws.Range("H2").Value = 120942.39513674159
$ws.Range("I2").Value = 118215.58250770152
$ws.Range("J2").Value = 111468.60834011578
$ws.Range("K2").Value = 108513.26254136374
$ws.Range("L2").Value = 100070.40199206585
$ws.Range("M2").Value = 97200.514426351001
$ws.Range("N2").Value = 94828.349484629987
$ws.Range("O2").Value = 92790.163275456362
$ws.Range("P2").Value = 90981.073576947689
$ws.Range("Q2").Value = 89457.401044596205
$ws.Range("R2").Value = 88233.458522771427
$ws.Range("S2").Value = 87279.701108456065
$ws.Range("T2").Value = 87639.912345720775
$ws.Range("U2").Value = 88921.979271717006
$ws.Range("V2").Value = 91019.505562573118
$ws.Range("W2").Value = 92072.174717311747
$ws.Range("X2").Value = 93431.448913084227
$ws.Range("Y2").Value = 95105.947297685983
$ws.Range("Z2").Value = 97105.736026253639
$ws.Range("AA2").Value = 99442.396211846062
$ws.Range("AB2").Value = 102129.07877060529
$ws.Range("AC2").Value = 105180.53591200689
$ws.Range("AD2").Value = 105109.33299626803
$ws.Range("AE2").Value = 106224.67976913461
$ws.Range("AF2").Value = 90162.96309238211

# Row 3 - Gale
$ws.Range("B3").Value = 146313.39580520417
$ws.Range("C3").Value = 146827.02811716936
$ws.Range("D3").Value = 148035.20963514346
$ws.Range("E3").Value = 137225.56905248243
$ws.Range("F3").Value = 133742.50143629144
$ws.Range("G3").Value = 123600.39784684683
$ws.Range("H3").Value = 121286.64585898256
$ws.Range("I3").Value = 118824.95991016011
$ws.Range("J3").Value = 114184.01144804618
$ws.Range("K3").Value = 114303.58014003787
$ws.Range("L3").Value = 109817.17055619013
$ws.Range("M3").Value = 108323.3258593293
$ws.Range("N3").Value = 107128.31842422836
$ws.Range("O3").Value = 106215.50069779737
$ws.Range("P3").Value = 105824.1190411592
$ws.Range("Q3").Value = 105897.06283442685
$ws.Range("R3").Value = 106353.04615057832
$ws.Range("S3").Value = 107163.05422865818
$ws.Range("T3").Value = 109371.72724991434
$ws.Range("U3").Value = 112591.3240874981
$ws.Range("V3").Value = 116716.32678416962
$ws.Range("W3").Value = 119888.75207185745
$ws.Range("X3").Value = 123461.65845937163
$ws.Range("Y3").Value = 127445.60792789987
$ws.Range("Z3").Value = 131849.30246469573
$ws.Range("AA3").Value = 136683.63772304164
$ws.Range("AB3").Value = 141957.97004292294
$ws.Range("AC3").Value = 147681.76018762818
$ws.Range("AD3").Value = 142044.15440705433
$ws.Range("AE3").Value = 139943.08492937451
$ws.Range("AF3").Value = 101697.54816029593

# Row 4 - Storm
$ws.Range("B4").Value = 146313.39580520417
$ws.Range("C4").Value = 146827.02811716936
$ws.Range("D4").Value = 148425.74962877488
$ws.Range("E4").Value = 138201.04137950164
$ws.Range("F4").Value = 135339.3915047746
$ws.Range("G4").Value = 125276.85616869353
$ws.Range("H4").Value = 123475.44645062975
$ws.Range("I4").Value = 121646.6707742135
$ws.Range("J4").Value = 119194.73116193741
$ws.Range("K4").Value = 122408.71465070732
$ws.Range("L4").Value = 121600.24121365207
$ws.Range("M4").Value = 122111.70953142329
$ws.Range("N4").Value = 123098.77111106152
$ws.Range("O4").Value = 124549.79040905094
$ws.Range("P4").Value = 126471.68999137332
$ws.Range("Q4").Value = 128874.40585720228
$ws.Range("R4").Value = 131820.27151360729
$ws.Range("S4").Value = 135290.58393349499
$ws.Range("T4").Value = 140340.04897952243
$ws.Range("U4").Value = 146592.85787902056
$ws.Range("V4").Value = 153961.40670160003
$ws.Range("W4").Value = 160599.56091487332
$ws.Range("X4").Value = 167880.70752498737
$ws.Range("Y4").Value = 175831.87187214554
$ws.Range("Z4").Value = 184481.53452999424
$ws.Range("AA4").Value = 193857.56976193868
$ws.Range("AB4").Value = 203986.0562376813
$ws.Range("AC4").Value = 214889.77639478218
$ws.Range("AD4").Value = 201580.05265164858
$ws.Range("AE4").Value = 195111.73962818645
$ws.Range("AF4").Value = 120836.92506986743
